$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing date strings in rows 2 and 3 (column A) to the new dithered values.
$ws.Range("A2").Value = "Patient X has an appointment on 1/2/2005 08:52 PM PST and another at 12/3/2006 09:53 PM PDT."
$ws.Range("A3").Value = "Can we see Patient X at 3/14/2007 11:00 AM PST and also at 04/15/2008 03:53 PM PDT?"

# Add two new test rows exercising month/day-only and month-year-only patterns.
$ws.Range("A7").Value = "Let's try just month/day without year like 11/28 or 3/2 and so forth."
$ws.Range("B7").Value = "ymuw64y7mu3w4e6"

$ws.Range("A8").Value = "What about just June 1999?"
$ws.Range("B8").Value = "m75ir67i ,r7i68oi"

# Move the active selection to match the saved workbook state.
$ws.Range("A12").Select()
